{"js": "// The document contains a closing line, italicized, that reads:\n//   \"Copia Verdadera del Original\"\n// It must be replaced with the merge-field placeholder:\n//   \"{o1}\"\n// (the template's \"original document checked/unchecked\" token), while\n// keeping the run's existing formatting (italic, theme fonts, language).\n\nconst body = context.document.body;\n\n// Search for the exact phrase so we only touch the one run that holds it.\nconst results = body.search(\"Copia Verdadera del Original\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Replacing (not inserting beside) keeps the original run's formatting.\n  results.items[i].insertText(\"{o1}\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a closing line, italicized, that reads:\n#   \"Copia Verdadera del Original\"\n# It must be replaced with the merge-field placeholder:\n#   \"{o1}\"\n# (the template's \"original document checked/unchecked\" token), while\n# keeping the run's existing formatting (italic, theme fonts, language).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Copia Verdadera del Original\"\n$find.Replacement.Text = \"{o1}\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
